# Add a new "23-nov" column (CM) to the right of the existing "22-nov" column (CL),
# mirroring the style of the CL column and filling in the new day's counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell (row 1): text value, same style as the other date headers (CL1)
$ws.Range("CM1").Value = "23-nov"
$ws.Range("CM1").NumberFormat = $ws.Range("CL1").NumberFormat

# Data cells (rows 2-11): numeric values, same style as the CL column
$values = @(13, 9, 8, 11, 11, 14, 11, 10, 17, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $srcCell = $ws.Cells.Item($row, 90)  # column 90 = CL
    $cell = $ws.Cells.Item($row, 91)     # column 91 = CM
    $cell.Value = $values[$i]
    $cell.HorizontalAlignment = $srcCell.HorizontalAlignment
    $cell.NumberFormat = $srcCell.NumberFormat
}

# Update the selection to match the new last cell, as seen in the diff
$ws.Range("CM11").Select()
